$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2755
$ws1.Range("F4").Value = 360
$ws1.Range("F5").Value = 1536
$ws1.Range("F6").Value = 1148
$ws1.Range("F11").Value = 552
$ws1.Range("F12").Value = 9334
$ws1.Range("F14").Value = 2504
$ws1.Range("F15").Value = 9
$ws1.Range("F19").Value = 644
$ws1.Range("F21").Value = 1188
$ws1.Range("F23").Value = 2922
$ws1.Range("F24").Value = 2212
$ws1.Range("F25").Value = 1915
$ws1.Range("F29").Value = 1551
$ws1.Range("F30").Value = 290
$ws1.Range("F31").Value = 5
$ws1.Range("F32").Value = 163
$ws1.Range("F33").Value = 217
$ws1.Range("F34").Value = 28
$ws1.Range("F37").Value = 302
$ws1.Range("F38").Value = 498
$ws1.Range("F39").Value = 17
$ws1.Range("F40").Value = 81
$ws1.Range("F41").Value = 970
$ws1.Range("F42").Value = 77
$ws1.Range("F43").Value = 1428
$ws1.Range("F44").Value = 8
$ws1.Range("F45").Value = 315
$ws1.Range("F46").Value = 14
$ws1.Range("F47").Value = 191
$ws1.Range("F48").Value = 677
$ws1.Range("F49").Value = 304

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 23
$ws2.Range("F5").Value = 25

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2755
$ws4.Range("F3").Value = 360
$ws4.Range("F4").Value = 1536
$ws4.Range("F6").Value = 1148
$ws4.Range("F9").Value = 9334
$ws4.Range("F11").Value = 2504
$ws4.Range("F12").Value = 23
$ws4.Range("F13").Value = 9
$ws4.Range("F17").Value = 644
$ws4.Range("F18").Value = 1189
$ws4.Range("F20").Value = 2212
$ws4.Range("F21").Value = 1915
$ws4.Range("F23").Value = 1551
$ws4.Range("F24").Value = 290
$ws4.Range("F25").Value = 5
$ws4.Range("F26").Value = 163
$ws4.Range("F27").Value = 217
$ws4.Range("F28").Value = 28
$ws4.Range("F31").Value = 302
$ws4.Range("F32").Value = 498
$ws4.Range("F33").Value = 25
$ws4.Range("F36").Value = 17
$ws4.Range("F37").Value = 81
$ws4.Range("F38").Value = 971
$ws4.Range("F40").Value = 77
$ws4.Range("F41").Value = 1428
$ws4.Range("F42").Value = 8
$ws4.Range("F44").Value = 315
$ws4.Range("F45").Value = 14
$ws4.Range("F46").Value = 191
$ws4.Range("F47").Value = 677
$ws4.Range("F48").Value = 304
